# Append a new submission row (row 18) to the bottom of the report table,
# mirroring the other form-response rows already on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Column C holds a quantity that happens to be the literal text "2" (like the
# other rows in this sheet, e.g. C14/C17). Force text formatting first so
# Excel doesn't silently coerce the numeric-looking string into a number.
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 2).Value = "أحمد شريم"
$ws.Cells.Item($row, 3).Value = "2"
$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 2"
$ws.Cells.Item($row, 6).Value = "C2"
$ws.Cells.Item($row, 7).Value = "IDRF"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٢٢:٣٤ م"
